$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Price (D) and Volume (E) columns so that
# numeric-looking strings (e.g. "23.90") are stored as text, matching the
# original inlineStr cell type, not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.585.60'
$ws.Range("E2").Value = '  -2.48%  '
$ws.Range("D3").Value = '1.664.10'
$ws.Range("E3").Value = '  -3.67%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '215.07'
$ws.Range("E5").Value = '  -1.98%  '
$ws.Range("E6").Value = '  -2.21%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("D8").Value = '23.90'
$ws.Range("E8").Value = '  -2.06%  '
$ws.Range("E9").Value = '  -0.77%  '
$ws.Range("E10").Value = '  -2.12%  '
$ws.Range("E11").Value = '  -2.17%  '
$ws.Range("D12").Value = '1.900.73'
$ws.Range("E12").Value = '  -3.60%  '
$ws.Range("D13").Value = '1.654.82'
$ws.Range("E13").Value = '  -4.23%  '
$ws.Range("E14").Value = '  -3.42%  '
$ws.Range("E15").Value = '  -0.64%  '
$ws.Range("D16").Value = '66.43'
$ws.Range("E16").Value = '  -1.89%  '
$ws.Range("D17").Value = '27.579.89'
$ws.Range("E17").Value = '  -2.31%  '
$ws.Range("D18").Value = '241.81'
$ws.Range("E18").Value = '  -1.06%  '
$ws.Range("D19").Value = '0.0₃0729'
$ws.Range("E19").Value = '  -3.59%  '
$ws.Range("E20").Value = '  -4.77%  '
$ws.Range("E21").Value = '  +0.20%  '
$ws.Range("E22").Value = '  -3.43%  '
$ws.Range("D23").Value = '9.34'
$ws.Range("E23").Value = '  -3.67%  '
$ws.Range("E24").Value = '  -3.68%  '
$ws.Range("D25").Value = '147.01'
$ws.Range("E25").Value = '  -1.53%  '
$ws.Range("D26").Value = '7.22'
$ws.Range("E26").Value = '  -4.13%  '
$ws.Range("E27").Value = '  -2.00%  '
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("E29").Value = '  -2.29%  '
$ws.Range("E30").Value = '  +2.39%  '
$ws.Range("E31").Value = '  -2.05%  '
$ws.Range("E32").Value = '  -2.67%  '
$ws.Range("D33").Value = '1.461.28'
$ws.Range("E33").Value = '  -2.87%  '
$ws.Range("E34").Value = '  -4.99%  '
$ws.Range("E35").Value = '  -5.08%  '
$ws.Range("E36").Value = '  -1.23%  '
$ws.Range("D37").Value = '0.923'
$ws.Range("E37").Value = '  -4.45%  '
$ws.Range("E38").Value = '  -1.74%  '
$ws.Range("E39").Value = '  -5.82%  '
$ws.Range("D40").Value = '1.04'
$ws.Range("E40").Value = '  -2.93%  '
$ws.Range("D41").Value = '69.62'
$ws.Range("E41").Value = '  -2.04%  '
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '5.41'
$ws.Range("E43").Value = '  -6.37%  '
$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").Value = '2.22'
$ws.Range("E44").Value = '  -3.47%  '
$ws.Range("D45").Value = '1.808.11'
$ws.Range("E45").Value = '  -3.64%  '
$ws.Range("E46").Value = '  -2.37%  '
$ws.Range("E47").Value = '  -1.79%  '
$ws.Range("E48").Value = '  -2.14%  '
$ws.Range("E49").Value = '  -5.68%  '
$ws.Range("D50").Value = '0.102'
$ws.Range("E50").Value = '  -2.77%  '
$ws.Range("D51").Value = '7.87'
$ws.Range("E51").Value = '  -4.79%  '

# Restore the default (Normal) style so no stray style/format is left on
# the cells beyond what the diff expects.
$ws.Range("D2:E51").Style = "Normal"

Write-Host "Applied cryptos list update"
